# Time tracker update for week of 4-26 (also fills in notes for the
# weeks of 4-12 and 4-19 that were left blank).
#
# The document is a single table; each row is "Week beginning ..."
# followed by a few content paragraphs (mostly empty placeholders) and
# a trailing blank paragraph. We locate each week by its heading text
# and then fill in / replace the blank paragraph(s) that follow it.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParaIndex {
    # Returns the 1-based index of the first paragraph whose text
    # contains $needle, searching from $start (1-based, inclusive).
    param([string]$needle, [int]$start = 1)
    for ($i = $start; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        if ($t -like "*$needle*") { return $i }
    }
    return -1
}

function Set-ParaXml {
    # Replaces the contents of paragraph $idx with the supplied run(s)
    # of OOXML (one or more <w:p> elements as a string).
    param([int]$idx, [string]$xml)
    $p = $d.Paragraphs.Item($idx)
    $p.Range.InsertXML($xml)
}

function Make-TextParaXml {
    param([string]$text)
    return "<w:p $wNs><w:r><w:t>$text</w:t></w:r></w:p>"
}

# ---------------------------------------------------------------------
# Week beginning April 12th: the first blank paragraph after "Read
# documentation for 'beautiful soup' library, " gets two sentences
# (typed as two separate runs); the following blank paragraph stays
# untouched.
# ---------------------------------------------------------------------
$aprReadIdx = Find-ParaIndex "Read documentation for"
$apr12BlankIdx = $aprReadIdx + 1

$xmlApr12 = "<w:p $wNs>" +
    "<w:r><w:t>Wrote scraper for initial data</w:t></w:r>" +
    "<w:r><w:t>, mostly working with a few issues.</w:t></w:r>" +
    "</w:p>"
Set-ParaXml $apr12BlankIdx $xmlApr12

# ---------------------------------------------------------------------
# Week beginning April 19th: the two blank paragraphs right after the
# heading become three note paragraphs (net +1 paragraph); the last
# blank paragraph of the row is left alone.
# ---------------------------------------------------------------------
$apr19HeadingIdx = Find-ParaIndex "April 19th"
$apr19FirstBlankIdx = $apr19HeadingIdx + 1
$apr19SecondBlankIdx = $apr19HeadingIdx + 2

# Drop the second blank paragraph entirely; its content will be
# re-created (three paragraphs) on top of the first blank paragraph.
$d.Paragraphs.Item($apr19SecondBlankIdx).Range.Delete()

$xmlApr19 = (Make-TextParaXml "Finished scraper, successfully scrapes venue, artist, and show data.") +
    (Make-TextParaXml "Datetime is giving me an issue, read a ton about UTC time, databases and converting datetime objects.") +
    (Make-TextParaXml "Looked at a bug, notes page is showing time in UTC instead of local time.")
Set-ParaXml $apr19FirstBlankIdx $xmlApr19

# ---------------------------------------------------------------------
# Week beginning April 26th: the bold, empty placeholder paragraph
# right after the heading becomes a normal (non-bold) text paragraph;
# the blank paragraph that used to follow it is removed (net -1
# paragraph). The remaining trailing blank paragraph is left alone.
# ---------------------------------------------------------------------
$apr26HeadingIdx = Find-ParaIndex "April 26th"
$apr26FirstBlankIdx = $apr26HeadingIdx + 1
$apr26SecondBlankIdx = $apr26HeadingIdx + 2

$d.Paragraphs.Item($apr26SecondBlankIdx).Range.Delete()

$xmlApr26 = Make-TextParaXml "Worked on a failing unit test, decided the regex needed a pipe character, return was put in a list first before counting."
Set-ParaXml $apr26FirstBlankIdx $xmlApr26
